$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("General")

# Fix the responsible for "Mejorar las preguntas" (row 4): Adrián -> Adrián y Diego
$ws.Range("E4").Value = "Adrián y Diego"

# New row 5 entry: mark Status "ok" (matching the style used for the other "ok" status
# cells, e.g. C2) and set Responsible to "Diego e Iván"
$ws.Range("C5").Value = "ok"
$ws.Range("C5").Style = $ws.Range("C2").Style

$ws.Range("E5").Value = "Diego e Iván"

# Update the active selection to C4
$ws.Range("C4").Select()
